$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Gnai2", "Adcy1", "ECs", 3, 1, 211.4614666666667, 634.3844, 0.2421062275331183, 0.2421062275331183, 1, 0.3333333333333333, 0.02287433333333333, 0.068623, 0.08516968032358596, 0.08516968032358598, 4.837040075688889, 43.5333606812, 0.02062011000334505, 0.02062011000334505),
    @("ECs", "Gnai2", "Adcy1", "FAPs", 3, 1, 211.4614666666667, 634.3844, 0.2421062275331183, 0.2421062275331183, 3, 1, 0.1212106666666667, 0.363632, 0.4513125511188116, 0.4513125511188117, 25.63138534897778, 230.6824681408, 0.1092655791897231, 0.1092655791897231),
    @("ECs", "Gnai2", "Adcy1", "MuSCs", 3, 1, 211.4614666666667, 634.3844, 0.2421062275331183, 0.2421062275331183, 3, 1, 0.1244886666666667, 0.373466, 0.4635177685576024, 0.4635177685576025, 26.32455603671111, 236.9210043304, 0.1122205383400502, 0.1122205383400502),
    @("FAPs", "Gnai2", "Adcy1", "ECs", 3, 1, 77.232923, 231.698769, 0.08842543241393927, 0.08842543241393927, 1, 0.3333333333333333, 0.02287433333333333, 0.068623, 0.08516968032358596, 0.08516968032358598, 1.766651625009667, 15.899864625087, 0.007531165811170064, 0.007531165811170065),
    @("FAPs", "Gnai2", "Adcy1", "FAPs", 3, 1, 77.232923, 231.698769, 0.08842543241393927, 0.08842543241393927, 3, 1, 0.1212106666666667, 0.363632, 0.4513125511188116, 0.4513125511188117, 9.361454085445335, 84.25308676900801, 0.03990750748651899, 0.03990750748651899),
    @("FAPs", "Gnai2", "Adcy1", "MuSCs", 3, 1, 77.232923, 231.698769, 0.08842543241393927, 0.08842543241393927, 3, 1, 0.1244886666666667, 0.373466, 0.4635177685576024, 0.4635177685576025, 9.614623607039332, 86.53161246335399, 0.04098675911625022, 0.04098675911625022),
    @("Inflammatory-Mac", "Gnai2", "Adcy1", "ECs", 3, 1, 174.3107043333333, 522.932113, 0.199571617988009, 0.199571617988009, 1, 0.3333333333333333, 0.02287433333333333, 0.068623, 0.08516968032358596, 0.08516968032358598, 3.987241154488777, 35.885170390399, 0.01699745090569955, 0.01699745090569955),
    @("Inflammatory-Mac", "Gnai2", "Adcy1", "FAPs", 3, 1, 174.3107043333333, 522.932113, 0.199571617988009, 0.199571617988009, 3, 1, 0.1212106666666667, 0.363632, 0.4513125511188116, 0.4513125511188117, 21.12831667937956, 190.154850114416, 0.09006917604507728, 0.09006917604507728),
    @("Inflammatory-Mac", "Gnai2", "Adcy1", "MuSCs", 3, 1, 174.3107043333333, 522.932113, 0.199571617988009, 0.199571617988009, 3, 1, 0.1244886666666667, 0.373466, 0.4635177685576024, 0.4635177685576025, 21.69970716818422, 195.297364513658, 0.09250499103723221, 0.09250499103723221),
    @("MuSCs", "Gnai2", "Adcy1", "ECs", 3, 1, 28.53474833333333, 85.60424499999999, 0.03266997236655063, 0.03266997236655063, 1, 0.3333333333333333, 0.02287433333333333, 0.068623, 0.08516968032358596, 0.08516968032358598, 0.6527133449594443, 5.874420104635, 0.002782491102639504, 0.002782491102639505),
    @("MuSCs", "Gnai2", "Adcy1", "FAPs", 3, 1, 28.53474833333333, 85.60424499999999, 0.03266997236655063, 0.03266997236655063, 3, 1, 0.1212106666666667, 0.363632, 0.4513125511188116, 0.4513125511188117, 3.458715868648889, 31.12844281784, 0.01474436857372904, 0.01474436857372905),
    @("MuSCs", "Gnai2", "Adcy1", "MuSCs", 3, 1, 28.53474833333333, 85.60424499999999, 0.03266997236655063, 0.03266997236655063, 3, 1, 0.1244886666666667, 0.373466, 0.4635177685576024, 0.4635177685576025, 3.552252773685554, 31.97027496317, 0.01514311269018208, 0.01514311269018208),
    @("Neutrophils", "Gnai2", "Adcy1", "ECs", 3, 1, 230.32901, 690.98703, 0.263708034289011, 0.263708034289011, 1, 0.3333333333333333, 0.02287433333333333, 0.068623, 0.08516968032358596, 0.08516968032358598, 5.268622551076667, 47.41760295969, 0.02245992897915631, 0.02245992897915631),
    @("Neutrophils", "Gnai2", "Adcy1", "FAPs", 3, 1, 230.32901, 690.98703, 0.263708034289011, 0.263708034289011, 3, 1, 0.1212106666666667, 0.363632, 0.4513125511188116, 0.4513125511188117, 27.91833285477334, 251.26499569296, 0.1190147457055006, 0.1190147457055006),
    @("Neutrophils", "Gnai2", "Adcy1", "MuSCs", 3, 1, 230.32901, 690.98703, 0.263708034289011, 0.263708034289011, 3, 1, 0.1244886666666667, 0.373466, 0.4635177685576024, 0.4635177685576025, 28.67335134955333, 258.06016214598, 0.1222333596043541, 0.1222333596043541),
    @("Resolving-Mac", "Gnai2", "Adcy1", "ECs", 3, 1, 151.5554656666667, 454.666397, 0.1735187154093718, 0.1735187154093718, 1, 0.3333333333333333, 0.02287433333333333, 0.068623, 0.08516968032358596, 0.08516968032358598, 3.466730240147889, 31.200572161331, 0.01477853352157549, 0.01477853352157549),
    @("Resolving-Mac", "Gnai2", "Adcy1", "FAPs", 3, 1, 151.5554656666667, 454.666397, 0.1735187154093718, 0.1735187154093718, 3, 1, 0.1212106666666667, 0.363632, 0.4513125511188116, 0.4513125511188117, 18.37013903043378, 165.331251273904, 0.07831117411826265, 0.07831117411826265),
    @("Resolving-Mac", "Gnai2", "Adcy1", "MuSCs", 3, 1, 151.5554656666667, 454.666397, 0.1735187154093718, 0.1735187154093718, 3, 1, 0.1244886666666667, 0.373466, 0.4635177685576024, 0.4635177685576025, 18.86693784688911, 169.802440622002, 0.08042900776953368, 0.08042900776953368)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

Write-Host "Updated $($data.Count) rows"
